# Ambermoon Advanced workbook update:
# "Added ancient key and chests for manyeyes' castle 2"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Todo sheet: tidy up two finished/renamed todo items
# ---------------------------------------------------------------
$wsTodo = $wb.Worksheets.Item("Todo")
$wsTodo.Range("A3").Value = "Add gatekeeper chests"
$wsTodo.Range("A5").Value = "Finish manyeyes castle"

# ---------------------------------------------------------------
# Items sheet: clarify existing cursed weapon text and add the
# new ancient key item used to unlock the manyeyes' castle 2 boss room
# ---------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items")
$wsItems.Range("D14").Value = "Cursed weapon from manyeyes' castle 1"

$wsItems.Range("A15").Value = 416
$wsItems.Range("B15").Value = "Alter Schlüssel / Ancient Key"
$wsItems.Range("C15").Value = "Key"
$wsItems.Range("D15").Value = "Opens the boss room in manyeyes' castle 2"

# ---------------------------------------------------------------
# Chests sheet: add the three new chests found in manyeyes' castle 2
# ---------------------------------------------------------------
$wsChests = $wb.Worksheets.Item("Chests")
$wsChests.Range("A19").Value = 149
$wsChests.Range("B19").Value = "Manyeyes'c castle 2 (462)"
$wsChests.Range("C19").Value = "1x Ancient Key"

$wsChests.Range("A20").Value = 150
$wsChests.Range("B20").Value = "Manyeyes'c castle 2 (462)"
$wsChests.Range("C20").Value = "2x Healing Potion II, 2x Spell Potion III, 4x Antidot"

$wsChests.Range("A21").Value = 151
$wsChests.Range("B21").Value = "Manyeyes'c castle 2 (462)"
$wsChests.Range("C21").Value = "10x Healing Potion I, 5x Spell Potion I, 1x Firebrand, 150 Gold"

# ---------------------------------------------------------------
# Restore the view/selection state left behind by the edit session
# ---------------------------------------------------------------
$wsItems.Range("G14").Select()
$wsTodo.Range("E8").Select()

$wsChests.Activate()
$wsChests.Range("C22").Select()
